$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell holding the default (no explicit number format) style used on this sheet
$defaultStyleRef = $ws.Range("B2")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.240.06"
$ws.Range("D2").Style = $defaultStyleRef.Style
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.279.58"
$ws.Range("D3").Style = $defaultStyleRef.Style
$ws.Range("E3").Value = "  +2.64%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.83"
$ws.Range("D5").Style = $defaultStyleRef.Style
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.84"
$ws.Range("D6").Style = $defaultStyleRef.Style
$ws.Range("E6").Value = "  +7.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.590"
$ws.Range("D7").Style = $defaultStyleRef.Style
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("D9").Style = $defaultStyleRef.Style
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.06"
$ws.Range("D10").Style = $defaultStyleRef.Style
$ws.Range("E10").Value = "  +6.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = $defaultStyleRef.Style
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.92"
$ws.Range("D12").Style = $defaultStyleRef.Style
$ws.Range("E12").Value = "  +1.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").Style = $defaultStyleRef.Style
$ws.Range("E13").Value = "  +1.62%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.633.37"
$ws.Range("D14").Style = $defaultStyleRef.Style
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.883"
$ws.Range("D15").Style = $defaultStyleRef.Style
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.64"
$ws.Range("D16").Style = $defaultStyleRef.Style
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.290.47"
$ws.Range("D17").Style = $defaultStyleRef.Style
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.206.35"
$ws.Range("D18").Style = $defaultStyleRef.Style
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.12"
$ws.Range("D19").Style = $defaultStyleRef.Style
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0998"
$ws.Range("D20").Style = $defaultStyleRef.Style
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("D21").Style = $defaultStyleRef.Style
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.33"
$ws.Range("D22").Style = $defaultStyleRef.Style
$ws.Range("E22").Value = "  +1.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.22"
$ws.Range("D23").Style = $defaultStyleRef.Style
$ws.Range("E23").Value = "  +1.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.80"
$ws.Range("D24").Style = $defaultStyleRef.Style
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.21"
$ws.Range("D25").Style = $defaultStyleRef.Style
$ws.Range("E25").Value = "  +4.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = $defaultStyleRef.Style
$ws.Range("E26").Value = "  -0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.29"
$ws.Range("D27").Style = $defaultStyleRef.Style
$ws.Range("E27").Value = "  +2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.12"
$ws.Range("D28").Style = $defaultStyleRef.Style
$ws.Range("E28").Value = "  +14.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = $defaultStyleRef.Style
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.55"
$ws.Range("D30").Style = $defaultStyleRef.Style
$ws.Range("E30").Value = "  +4.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.04"
$ws.Range("D31").Style = $defaultStyleRef.Style
$ws.Range("E31").Value = "  +5.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0887"
$ws.Range("D32").Style = $defaultStyleRef.Style
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.54"
$ws.Range("D33").Style = $defaultStyleRef.Style
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -1.02%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.08"
$ws.Range("D35").Style = $defaultStyleRef.Style
$ws.Range("E35").Value = "  +3.53%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.27"
$ws.Range("D36").Style = $defaultStyleRef.Style
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("D37").Style = $defaultStyleRef.Style
$ws.Range("E37").Value = "  +13.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("D38").Style = $defaultStyleRef.Style
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.51"
$ws.Range("D39").Style = $defaultStyleRef.Style
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.95"
$ws.Range("D40").Style = $defaultStyleRef.Style
$ws.Range("E40").Value = "  +6.50%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0327"
$ws.Range("D41").Style = $defaultStyleRef.Style
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.46"
$ws.Range("D42").Style = $defaultStyleRef.Style
$ws.Range("E42").Value = "  +27.24%  "
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.769.42"
$ws.Range("D44").Style = $defaultStyleRef.Style
$ws.Range("E44").Value = "  -6.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.208"
$ws.Range("D45").Style = $defaultStyleRef.Style
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "85.82"
$ws.Range("D46").Style = $defaultStyleRef.Style
$ws.Range("E46").Value = "  -4.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.40"
$ws.Range("D47").Style = $defaultStyleRef.Style
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.89"
$ws.Range("D48").Style = $defaultStyleRef.Style
$ws.Range("E48").Value = "  +3.45%  "
$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "75.72"
$ws.Range("D49").Style = $defaultStyleRef.Style
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.69"
$ws.Range("D50").Style = $defaultStyleRef.Style
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "104.66"
$ws.Range("D51").Style = $defaultStyleRef.Style
$ws.Range("E51").Value = "  +3.20%  "
